# Update the two mail-merge label blocks (serah terima BN dp2nt15 dan dp3n30)
#
# Note: this runtime's Range.Find.Execute always searches the whole
# document story regardless of which Range/Cell/Paragraph it was invoked
# on, so every Find below is issued against $d.Content and relies on
# wdReplaceOne (the 11th arg = 1) to hit exactly one, deterministic,
# in-document-order occurrence. The two ambiguous numeric values (the
# "42" that appears twice, and the "55"/"58" pair that would collide if
# done in the wrong order) are therefore replaced in an order chosen so
# that each Find.Execute call lands on the correct occurrence:
#
#   1) i1  -> i19                 (unique)
#   2) BAYU AJI SANTOSO -> DIMAS ARIFIN   (unique)
#   3) 42  -> 41   (ReplaceOne hits the FIRST "42", i.e. block 1's SEPATU;
#                   block 2's SEPATU "42" is left untouched, as required)
#   4) 58  -> 59   (done BEFORE step 5, while "58" is still unique - it is
#                   block 2's NO. TOPI value)
#   5) 55  -> 58   (block 1's NO. TOPI value; now safe since the only "58"
#                   left at this point is the freshly-written one from
#                   step 4, which Find no longer needs to touch)
#   6) i2  -> i20                 (unique)
#   7) PUJIANTO -> MOH. IRWAN TOPAN       (unique)

$d = $word.ActiveDocument

$d.Content.Find.Execute("i1", $true, $false, $false, $false, $false, $true, 0, $false, "i19", 1)
$d.Content.Find.Execute("BAYU AJI SANTOSO", $true, $false, $false, $false, $false, $true, 0, $false, "DIMAS ARIFIN", 1)
$d.Content.Find.Execute("42", $true, $false, $false, $false, $false, $true, 0, $false, "41", 1)
$d.Content.Find.Execute("58", $true, $false, $false, $false, $false, $true, 0, $false, "59", 1)
$d.Content.Find.Execute("55", $true, $false, $false, $false, $false, $true, 0, $false, "58", 1)
$d.Content.Find.Execute("i2", $true, $false, $false, $false, $false, $true, 0, $false, "i20", 1)
$d.Content.Find.Execute("PUJIANTO", $true, $false, $false, $false, $false, $true, 0, $false, "MOH. IRWAN TOPAN", 1)
